$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08780144004532819
$ws.Range("H2").Value = -7.305088008990253
$ws.Range("I2").Value = 44.65302825260931
$ws.Range("G3").Value = 0.09202401271684942
$ws.Range("H3").Value = 39.23221947467672
$ws.Range("G4").Value = 0.009203472610279827
$ws.Range("H4").Value = -1.832427162465355
$ws.Range("G5").Value = 0.05296606007633486
$ws.Range("H5").Value = 495.5698461288206
$ws.Range("G6").Value = -0.2254735375273857
$ws.Range("H6").Value = -1.933966054407766
$ws.Range("G7").Value = -0.2606665125120894
$ws.Range("H7").Value = -4.318474521871623
$ws.Range("G8").Value = -0.3924973142448263
$ws.Range("H8").Value = -6.007047102363014
$ws.Range("G9").Value = -0.352611259055911
$ws.Range("H9").Value = 11.54519382092288
$ws.Range("G10").Value = -0.0292761444383459
$ws.Range("H10").Value = -280.6624270126746
$ws.Range("G11").Value = 0.04262529026098123
$ws.Range("H11").Value = 365.3158362719031
$ws.Range("G12").Value = 0.2178775189809439
$ws.Range("H12").Value = -4.098477585233849
$ws.Range("G13").Value = 0.286535152397906
$ws.Range("H13").Value = 8.807821028511535
$ws.Range("G14").Value = -0.01895751064732112
$ws.Range("H14").Value = -98.14958613873105
$ws.Range("G15").Value = 0.02728944632635594
$ws.Range("H15").Value = 35.18516361872027
$ws.Range("G16").Value = 0.1193543737835493
$ws.Range("H16").Value = 1.148541263879328
$ws.Range("G17").Value = 0.1795338430622062
$ws.Range("H17").Value = -17.95848582926734
$ws.Range("G18").Value = 0.04239690756981041
$ws.Range("H18").Value = -29.88376206018753
$ws.Range("G19").Value = 0.08773393159087975
$ws.Range("H19").Value = -2.614486303705878
$ws.Range("G20").Value = -0.1487736286292619
$ws.Range("H20").Value = -2.228153623497549
$ws.Range("G21").Value = -0.1800433447937272
$ws.Range("H21").Value = 9.901696196330363
$ws.Range("G22").Value = 0.03735406771233536
$ws.Range("H22").Value = -31.3195905201351
$ws.Range("G23").Value = 0.0618586224879247
$ws.Range("H23").Value = 51.4652364659373
$ws.Range("G24").Value = 0.1080587264919119
$ws.Range("H24").Value = -6.632083166709856
$ws.Range("G25").Value = 0.1541904037717162
$ws.Range("H25").Value = 1.390340882666623
$ws.Range("G26").Value = 0.04008186082195228
$ws.Range("H26").Value = -24.19762383530551
$ws.Range("G27").Value = 0.04769684511658644
$ws.Range("H27").Value = -5.488681579355728
$ws.Range("G28").Value = 0.1639230147903397
$ws.Range("H28").Value = 7.201533646178758
$ws.Range("G29").Value = 0.1949935796757985
$ws.Range("H29").Value = 14.229673308847
$ws.Range("G30").Value = 0.02789225223439786
$ws.Range("H30").Value = 42.55406516175641
$ws.Range("G31").Value = 0.04379583579400411
$ws.Range("H31").Value = 351.2722149207809
$ws.Range("G32").Value = 0.03625063120518941
$ws.Range("H32").Value = -2.795788700318291
$ws.Range("G33").Value = 0.0204554949812015
$ws.Range("H33").Value = -21.63738135517511
$ws.Range("G34").Value = 0.07814092310363718
$ws.Range("H34").Value = -38.93714222735387
$ws.Range("G35").Value = 0.1451589308708243
$ws.Range("H35").Value = 12.8231476532991
$ws.Range("G36").Value = -0.0004743575117232557
$ws.Range("H36").Value = -103.155674980042
$ws.Range("G37").Value = 0.04670055068784153
$ws.Range("H37").Value = 204.9435897532105
$ws.Range("G38").Value = -0.00740678114589158
$ws.Range("H38").Value = -262.6484463483891
$ws.Range("G39").Value = 0.02863559196215204
$ws.Range("H39").Value = 185.7116642599964
$ws.Range("G40").Value = 0.1307558692232715
$ws.Range("H40").Value = -11.38138770856821
$ws.Range("G41").Value = 0.1670133006373962
$ws.Range("H41").Value = 3.478390424263764
$ws.Range("G42").Value = 0.05712942873300549
$ws.Range("H42").Value = -11.51628048846844
$ws.Range("G43").Value = 0.05021808988384746
$ws.Range("H43").Value = 44.46874002629264
$ws.Range("G44").Value = 0.0373975350966762
$ws.Range("H44").Value = 164.9905123105917
$ws.Range("G45").Value = 0.04062660362704599
$ws.Range("H45").Value = -1.050792665530195
$ws.Range("G46").Value = -0.04011426250070976
$ws.Range("H46").Value = 39.05447374819153
$ws.Range("G47").Value = -0.05160746553427206
$ws.Range("H47").Value = -24.92819367653398
$ws.Range("G48").Value = -0.08256421516848483
$ws.Range("H48").Value = 34.46005411156436
$ws.Range("G49").Value = -0.155711390964151
$ws.Range("H49").Value = 21.15136626984043
$ws.Range("G50").Value = 0.1326565210207497
$ws.Range("H50").Value = 21.84447176313874
$ws.Range("G51").Value = 0.1378309156399299
$ws.Range("H51").Value = 37.45863870133678
$ws.Range("G52").Value = 0.0536808935918659
$ws.Range("H52").Value = -9.960171038253637
$ws.Range("G53").Value = 0.07384411119627586
$ws.Range("H53").Value = 9.329553862048247
$ws.Range("G54").Value = -0.0455602742494054
$ws.Range("H54").Value = 34.8399655230363
$ws.Range("G55").Value = -0.04701755511397501
$ws.Range("H55").Value = 39.12497825148908
$ws.Range("G56").Value = 0.04470127403729982
$ws.Range("H56").Value = -2.454496904906267
$ws.Range("G57").Value = 0.06861820902896697
$ws.Range("H57").Value = 1227.166769896331
